$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-ParaText($para, $text) {
    # PowerPoint's minimal-diff text setter only rewrites the differing
    # substring when old/new text share a common prefix/suffix, which can
    # split a single run into several. Routing through a throwaway value
    # first forces a full single-run rewrite that keeps formatting (rPr).
    $para.Text = "X"
    $para.Text = $text
}

# --- TextBox 16: RACM intro paragraph ---
$tb16 = $s.Shapes.Item("TextBox 16")
Set-ParaText $tb16.TextFrame.TextRange "RACM is the Resource Access Control Management component of SciServer (http://www.sciserver.org), a system for advanced analysis in a collaborative environment for large scientific datasets."

# --- TextBox 17: RACM data model paragraph + shrink box height ---
$tb17 = $s.Shapes.Item("TextBox 17")
Set-ParaText $tb17.TextFrame.TextRange "RACM uses a flexible data model for representing who is allowed to do which actions on each resource."
$tb17.Height = 104.2078

# --- TextBox 20: Core Concepts bullets ---
$tb20 = $s.Shapes.Item("TextBox 20")
$tb20tr = $tb20.TextFrame.TextRange
Set-ParaText $tb20tr.Paragraphs(2,1) "Define and manage users and groups"
Set-ParaText $tb20tr.Paragraphs(3,1) "Define and manage system resources"
Set-ParaText $tb20tr.Paragraphs(4,1) "Define and manage access controls between users and resources"
Set-ParaText $tb20tr.Paragraphs(6,1) "Support nested groups"
Set-ParaText $tb20tr.Paragraphs(7,1) "Coming Soon: regions and domains"

# --- TextBox 21: Major Features bullets ---
$tb21 = $s.Shapes.Item("TextBox 21")
$tb21tr = $tb21.TextFrame.TextRange
Set-ParaText $tb21tr.Paragraphs(4,1) "Workspaces for intuitive user access"

# --- TextBox 22: Flexibility and Application bullets ---
$tb22 = $s.Shapes.Item("TextBox 22")
$tb22tr = $tb22.TextFrame.TextRange
$tb22OrigHeight = $tb22.Height
Set-ParaText $tb22tr.Paragraphs(3,1) "Abstract model allows new resource types to be defined"
Set-ParaText $tb22tr.Paragraphs(4,1) "Can be used by external applications as an access control framework in itself"
# Re-applying the original autofit height avoids an off-by-one-EMU drift
# that the autofit recalculation introduces when runs are rewritten.
$tb22.Height = $tb22OrigHeight

# --- TextBox 52: User Workspaces bullets (split bullet 1 into two, merge bullet 2's runs) ---
$tb52 = $s.Shapes.Item("TextBox 52")
$tb52tr = $tb52.TextFrame.TextRange
$tb52tr.Text = "User created and managed" + [char]13 + "Easy to add resources and users" + [char]13 + "Shared workspace folder and database" + [char]13 + "Configurable permissions"
